$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1800.6666
$ws.Cells.Item(18, 9).Value = 1660.8889
$ws.Cells.Item(18, 10).Value = 2220
$ws.Cells.Item(18, 11).Value = 1660.8889
$ws.Cells.Item(18, 12).Value = 2220
$ws.Cells.Item(18, 13).Value = -1376.8889
$ws.Cells.Item(18, 14).Value = -2788

$ws.Cells.Item(62, 8).Value = 3007.2778
$ws.Cells.Item(62, 9).Value = 2469.0715
$ws.Cells.Item(62, 11).Value = 2469.0715
$ws.Cells.Item(62, 13).Value = -1845.0715

$ws.Cells.Item(64, 8).Value = 4845.1113
$ws.Cells.Item(64, 10).Value = 4754.5
$ws.Cells.Item(64, 12).Value = 4754.5
$ws.Cells.Item(64, 14).Value = -5250.5

$ws.Cells.Item(65, 8).Value = 3007.2778
$ws.Cells.Item(65, 9).Value = 2469.0715
$ws.Cells.Item(65, 11).Value = 12345.3575
$ws.Cells.Item(65, 13).Value = -9225.3575

$ws.Cells.Item(67, 8).Value = 4845.1113
$ws.Cells.Item(67, 10).Value = 4754.5
$ws.Cells.Item(67, 12).Value = 4754.5
$ws.Cells.Item(67, 14).Value = -6470.5

$ws.Cells.Item(100, 8).Value = 5047.722
$ws.Cells.Item(100, 9).Value = 3534.4285
$ws.Cells.Item(100, 10).Value = 6010.727
$ws.Cells.Item(100, 11).Value = 3534.4285
$ws.Cells.Item(100, 12).Value = 6010.727
$ws.Cells.Item(100, 13).Value = -2993.4285
$ws.Cells.Item(100, 14).Value = -7092.727

$ws.Cells.Item(107, 8).Value = 1470.7894
$ws.Cells.Item(107, 9).Value = 1219.9166
$ws.Cells.Item(107, 10).Value = 1900.8572
$ws.Cells.Item(107, 11).Value = 1219.9166
$ws.Cells.Item(107, 12).Value = 1900.8572
$ws.Cells.Item(107, 13).Value = 700.0834
$ws.Cells.Item(107, 14).Value = -5740.8572

$ws.Cells.Item(132, 8).Value = 4674.5884
$ws.Cells.Item(132, 9).Value = 2406.1738
$ws.Cells.Item(132, 11).Value = 7218.5214
$ws.Cells.Item(132, 13).Value = -4688.5214

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1910.5209
$ws.Cells.Item(32, 9).Value = 1697.841
$ws.Cells.Item(32, 11).Value = 1697.841
$ws.Cells.Item(32, 13).Value = -1410.841

$ws.Cells.Item(102, 8).Value = 2122.6
$ws.Cells.Item(102, 9).Value = 1961.2307
$ws.Cells.Item(102, 11).Value = 1961.2307
$ws.Cells.Item(102, 13).Value = -339.2307000000001

$ws.Cells.Item(132, 8).Value = 3576187.5
$ws.Cells.Item(132, 9).Value = 4935.185
$ws.Cells.Item(132, 11).Value = 14805.555
$ws.Cells.Item(132, 13).Value = -12275.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1330.3721
$ws.Cells.Item(94, 9).Value = 988.3333
$ws.Cells.Item(94, 11).Value = 988.3333
$ws.Cells.Item(94, 13).Value = -537.3333

$ws.Cells.Item(99, 8).Value = 3027.3076
$ws.Cells.Item(99, 9).Value = 3195.6
$ws.Cells.Item(99, 10).Value = 2466.3333
$ws.Cells.Item(99, 11).Value = 3195.6
$ws.Cells.Item(99, 12).Value = 2466.3333
$ws.Cells.Item(99, 13).Value = -1697.6
$ws.Cells.Item(99, 14).Value = -5462.3333

$ws.Cells.Item(105, 8).Value = 563991.2
$ws.Cells.Item(105, 9).Value = 806575
$ws.Cells.Item(105, 10).Value = 9514
$ws.Cells.Item(105, 11).Value = 806575
$ws.Cells.Item(105, 12).Value = 9514
$ws.Cells.Item(105, 13).Value = -804828
$ws.Cells.Item(105, 14).Value = -13008

$ws.Cells.Item(134, 8).Value = 20002640
$ws.Cells.Item(134, 9).Value = 3066.3333
$ws.Cells.Item(134, 10).Value = 50002000
$ws.Cells.Item(134, 11).Value = 9198.999899999999
$ws.Cells.Item(134, 12).Value = 150006000
$ws.Cells.Item(134, 13).Value = -6663.999899999999
$ws.Cells.Item(134, 14).Value = -150011070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1488.9362
$ws.Cells.Item(107, 9).Value = 1268.8055
$ws.Cells.Item(107, 11).Value = 1268.8055
$ws.Cells.Item(107, 13).Value = 651.1945000000001

$ws.Cells.Item(121, 8).Value = 73833
$ws.Cells.Item(121, 10).Value = 73833
$ws.Cells.Item(121, 12).Value = 73833
$ws.Cells.Item(121, 14).Value = -76453

$ws.Cells.Item(122, 8).Value = 4183.3335
$ws.Cells.Item(122, 9).Value = 4150
$ws.Cells.Item(122, 11).Value = 12450
$ws.Cells.Item(122, 13).Value = -10000

$ws.Cells.Item(132, 8).Value = 3175.9375
$ws.Cells.Item(132, 9).Value = 3100.0715
$ws.Cells.Item(132, 10).Value = 3707
$ws.Cells.Item(132, 11).Value = 9300.2145
$ws.Cells.Item(132, 12).Value = 11121
$ws.Cells.Item(132, 13).Value = -6770.2145
$ws.Cells.Item(132, 14).Value = -16181

$ws.Cells.Item(134, 8).Value = 3142.7334
$ws.Cells.Item(134, 9).Value = 3142.7334
$ws.Cells.Item(134, 11).Value = 9428.200199999999
$ws.Cells.Item(134, 13).Value = -6893.200199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 15163.047
$ws.Cells.Item(56, 9).Value = 15163.047
$ws.Cells.Item(56, 11).Value = 15163.047
$ws.Cells.Item(56, 13).Value = -14633.047

$ws.Cells.Item(57, 8).Value = 13891.167
$ws.Cells.Item(57, 10).Value = 19777.334
$ws.Cells.Item(57, 12).Value = 59332.00199999999
$ws.Cells.Item(57, 14).Value = -60450.00199999999

$ws.Cells.Item(115, 8).Value = 2864.5833
$ws.Cells.Item(115, 9).Value = 131.14285
$ws.Cells.Item(115, 11).Value = 393.42855
$ws.Cells.Item(115, 13).Value = 781.5714499999999

$ws.Cells.Item(122, 8).Value = 33352.2
$ws.Cells.Item(122, 9).Value = 66307.60000000001
$ws.Cells.Item(122, 11).Value = 596768.4
$ws.Cells.Item(122, 13).Value = -594318.4

$ws.Cells.Item(131, 8).Value = 4048.5588
$ws.Cells.Item(131, 9).Value = 2650.5833
$ws.Cells.Item(131, 11).Value = 7951.749899999999
$ws.Cells.Item(131, 13).Value = -2911.749899999999

$ws.Cells.Item(132, 8).Value = 3101
$ws.Cells.Item(132, 9).Value = 3004
$ws.Cells.Item(132, 10).Value = 3149.5
$ws.Cells.Item(132, 11).Value = 27036
$ws.Cells.Item(132, 12).Value = 28345.5
$ws.Cells.Item(132, 13).Value = -24506
$ws.Cells.Item(132, 14).Value = -33405.5

$ws.Cells.Item(137, 8).Value = 10648.272
$ws.Cells.Item(137, 9).Value = 3760
$ws.Cells.Item(137, 10).Value = 16388.5
$ws.Cells.Item(137, 11).Value = 11280
$ws.Cells.Item(137, 12).Value = 49165.5
$ws.Cells.Item(137, 13).Value = -6180
$ws.Cells.Item(137, 14).Value = -59365.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5639.4185
$ws.Cells.Item(70, 9).Value = 6030.25
$ws.Cells.Item(70, 11).Value = 6030.25
$ws.Cells.Item(70, 13).Value = -5760.25

$ws.Cells.Item(73, 8).Value = 5639.4185
$ws.Cells.Item(73, 9).Value = 6030.25
$ws.Cells.Item(73, 11).Value = 6030.25
$ws.Cells.Item(73, 13).Value = -5094.25

$ws.Cells.Item(97, 8).Value = 651
$ws.Cells.Item(97, 9).Value = 548.2857
$ws.Cells.Item(97, 10).Value = 770.8333
$ws.Cells.Item(97, 11).Value = 548.2857
$ws.Cells.Item(97, 12).Value = 770.8333
$ws.Cells.Item(97, 13).Value = -52.28570000000002
$ws.Cells.Item(97, 14).Value = -1762.8333

$ws.Cells.Item(102, 8).Value = 2855.875
$ws.Cells.Item(102, 9).Value = 2978.1428
$ws.Cells.Item(102, 10).Value = 2000
$ws.Cells.Item(102, 11).Value = 2978.1428
$ws.Cells.Item(102, 12).Value = 2000
$ws.Cells.Item(102, 13).Value = -1356.1428
$ws.Cells.Item(102, 14).Value = -5244

$ws.Cells.Item(132, 8).Value = 6588060
$ws.Cells.Item(132, 9).Value = 3672.0715
$ws.Cells.Item(132, 10).Value = 12733488
$ws.Cells.Item(132, 11).Value = 11016.2145
$ws.Cells.Item(132, 12).Value = 38200464
$ws.Cells.Item(132, 13).Value = -8486.2145
$ws.Cells.Item(132, 14).Value = -38205524

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4875.8
$ws.Cells.Item(61, 9).Value = 1677.4166
$ws.Cells.Item(61, 10).Value = 17669.334
$ws.Cells.Item(61, 11).Value = 1677.4166
$ws.Cells.Item(61, 12).Value = 17669.334
$ws.Cells.Item(61, 13).Value = -1475.4166
$ws.Cells.Item(61, 14).Value = -18073.334

$ws.Cells.Item(113, 8).Value = 4875.8
$ws.Cells.Item(113, 9).Value = 1677.4166
$ws.Cells.Item(113, 10).Value = 17669.334
$ws.Cells.Item(113, 11).Value = 1677.4166
$ws.Cells.Item(113, 12).Value = 17669.334
$ws.Cells.Item(113, 13).Value = 492.5834
$ws.Cells.Item(113, 14).Value = -22009.334

$ws.Cells.Item(122, 8).Value = 3545.6
$ws.Cells.Item(122, 9).Value = 3354.4546
$ws.Cells.Item(122, 11).Value = 10063.3638
$ws.Cells.Item(122, 13).Value = -7613.363799999999

$ws.Cells.Item(132, 8).Value = 3805
$ws.Cells.Item(132, 9).Value = 2357.6365
$ws.Cells.Item(132, 11).Value = 7072.9095
$ws.Cells.Item(132, 13).Value = -4542.9095

$ws.Cells.Item(136, 8).Value = 5401.2
$ws.Cells.Item(136, 9).Value = 4251.5
$ws.Cells.Item(136, 11).Value = 12754.5
$ws.Cells.Item(136, 13).Value = -10204.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 834552.2
$ws.Cells.Item(132, 9).Value = 1182.3334
$ws.Cells.Item(132, 11).Value = 3547.0002
$ws.Cells.Item(132, 13).Value = -1017.0002
